# Add compounds metrics calculation
# Fills rows 2-5 on Sheet1 with tool-quality metrics data (columns A-H).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @("aureme",       154, 491,  3915, 0.2387596899224806, 0.03784713688867043, 0.06533729316928298, 0.9685521748008985),
    @("merlin_blast",  164, 218,  3910, 0.4293193717277487, 0.04025527736867943, 0.07360861759425494, 0.9622641509433962),
    @("carveme",       411, 1172, 3654, 0.2596336070751737, 0.1011070110701107,  0.1455382436260623,  0.9314771590530176),
    @("kbase",         200, 338,  3868, 0.3717472118959108, 0.04916420845624386, 0.08684324793747285, 0.9559859154929577)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $col = $c + 1
        $ws.Cells.Item($row, $col).Value = $values[$c]
    }
}
